$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 88.69904661368808
$ws.Range("C2").Value = 96.97591407789436
$ws.Range("D2").Value = 98.71218898370287
$ws.Range("E2").Value = 98.59262851751754
$ws.Range("F2").Value = 98.35568571259847
$ws.Range("G2").Value = 97.56738495387086
$ws.Range("H2").Value = 97.32725144807199
$ws.Range("I2").Value = 96.13732136780457
$ws.Range("B3").Value = 81.90962221409752
$ws.Range("C3").Value = 97.7827893731686
$ws.Range("D3").Value = 99.20575852297434
$ws.Range("E3").Value = 98.67835946139346
$ws.Range("F3").Value = 98.53736482272471
$ws.Range("G3").Value = 97.91486804929083
$ws.Range("H3").Value = 97.53002160587134
$ws.Range("I3").Value = 95.83391430103309
$ws.Range("B4").Value = 85.36702098740928
$ws.Range("C4").Value = 96.13968860130888
$ws.Range("D4").Value = 99.01251692958789
$ws.Range("E4").Value = 98.40882978140412
$ws.Range("F4").Value = 98.46778682113062
$ws.Range("G4").Value = 97.60273362047018
$ws.Range("H4").Value = 97.70950222766976
$ws.Range("I4").Value = 95.95859431478002
$ws.Range("B5").Value = 87.24192865880424
$ws.Range("C5").Value = 95.61764349722598
$ws.Range("D5").Value = 98.5590893269721
$ws.Range("E5").Value = 98.52149474957137
$ws.Range("F5").Value = 98.35283057838525
$ws.Range("G5").Value = 97.93711659961268
$ws.Range("H5").Value = 97.4541398737331
$ws.Range("I5").Value = 96.1240712166569
$ws.Range("B6").Value = 87.35820418963948
$ws.Range("C6").Value = 96.62877182165764
$ws.Range("D6").Value = 99.18456649274428
$ws.Range("E6").Value = 98.83475658541579
$ws.Range("F6").Value = 98.39246238833009
$ws.Range("G6").Value = 98.14354718080673
$ws.Range("H6").Value = 97.34442529919663
$ws.Range("I6").Value = 96.34407731302136
